# Update CircadiPy cosinor results (sawtooth_10) after re-running the
# CircaDB / CircadiPy simulation analysis.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("E2").Value = 25.55000000000055
$ws.Range("G2").Value = 0.0011067150142271
$ws.Range("H2").Value = 0.002613270646225289
$ws.Range("K2").Value = 5.328976076040837
$ws.Range("L2").Value = "[2.2086512459124688, 8.449300906169205]"
$ws.Range("M2").Value = 0.000914818006285234
$ws.Range("N2").Value = 0.000914818006285234
$ws.Range("O2").Value = -1.522052897234695
$ws.Range("P2").Value = "[-2.314526719761851, -0.729579074707539]"
$ws.Range("Q2").Value = 0.0002034999909019675
$ws.Range("R2").Value = 0.0002034999909019675
$ws.Range("S2").Value = 10.66146766849726
$ws.Range("T2").Value = "[8.679978977780847, 12.642956359213665]"
$ws.Range("W2").Value = 6.189289289289423
$ws.Range("X2").Value = 2.966766766766827
$ws.Range("Y2").Value = 9.411811811812019

# ---- Row 3 ----
$ws.Range("E3").Value = 22.73000000000011
$ws.Range("G3").Value = [double]"3.65295153237355e-07"
$ws.Range("H3").Value = [double]"5.968065596189554e-06"
# I3 previously held a p_reject value; it is now blank (no longer rejected)
$ws.Range("I3").Value = ""
$ws.Range("K3").Value = 6.728538253900272
$ws.Range("L3").Value = "[3.918689234250836, 9.538387273549708]"
$ws.Range("M3").Value = [double]"4.081246445908349e-06"
$ws.Range("N3").Value = [double]"8.162492891816697e-06"
$ws.Range("O3").Value = 1.83023716155081
$ws.Range("P3").Value = "[1.3270791789938867, 2.333395144107734]"
$ws.Range("Q3").Value = [double]"9.722000982037571e-12"
$ws.Range("R3").Value = [double]"1.944400196407514e-11"
$ws.Range("S3").Value = 10.42647502345576
$ws.Range("T3").Value = "[8.76351636906151, 12.089433677850007]"
$ws.Range("W3").Value = 16.10894894894903
$ws.Range("X3").Value = 14.2887287287288
$ws.Range("Y3").Value = 17.92916916916926
